$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DividendHistory")

# A new dividend record (XD Date 30/12/2024, Pay Date 30/12/2024,
# Gross Dividend 0.008) needs to be inserted right under the header row,
# pushing every existing record down by one row.
#
# Row 3 (XD Date 27/06/2024) already carries a Gross Dividend of "0.008" -
# the same text the new row needs - and uses the sheet's normal (default)
# cell style. Copying that row down for the insert means the new row's
# Gross Dividend cell is already correct and keeps the plain-text storage
# / default style every other cell in the sheet uses, instead of having
# Excel re-interpret a freshly typed "0.008" as a number.
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(2).Insert()

# Only the two date cells need to change; the copied Gross Dividend cell
# (already "0.008") is left as-is.
$ws.Cells.Item(2, 1).Value = "30/12/2024"
$ws.Cells.Item(2, 2).Value = "30/12/2024"
